# Commit: Tue, Jun 23, 2020  4:05:06 AM
#
# Two changes:
#   1. Slide 16's table is re-styled with a different built-in PowerPoint
#      table style (tableStyleId GUID change).
#   2. The presentation's live theme (ppt/theme/theme2.xml, the part that
#      both the slide master and the presentation itself point at) is
#      switched from the custom "Integral" colour palette back to the
#      stock "Office Theme" colour palette. (Font scheme / format scheme
#      are identical between the two themes, so only the 12 theme colours
#      need to change.)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table on slide 16: apply the new built-in table style.
# ---------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{534DF2FA-408C-408A-9154-1780BE4904DA}")
    }
}

# ---------------------------------------------------------------------
# 2) Restore the stock "Office Theme" colour scheme on the deck's theme.
# ---------------------------------------------------------------------
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Index order per MsoThemeColorSchemeIndex:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hyperlink, 12 followed hyperlink
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeRGB[$i - 1]
}
